$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-18 (B column topic-association strings)
$ws.Range("B2").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.523|x10:0.000|x11:0.477|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B3").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:1.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B4").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:1.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B5").Value = "x1:0.000|x2:0.000|x3:0.164|x4:0.000|x5:0.000|x6:0.836|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B6").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:1.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B7").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:1.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B8").Value = "x1:0.000|x2:0.000|x3:1.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B9").Value = "x1:0.000|x2:0.000|x3:0.000|x4:1.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B10").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:1.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B11").Value = "x1:0.784|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.216|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B12").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:1.000|x16:0.000|x17:0.000"
$ws.Range("B13").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:1.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B14").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:1.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B15").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.760|x10:0.240|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B16").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.211|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.789|x17:0.000"
$ws.Range("B17").Value = "x1:0.000|x2:1.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B18").Value = "x1:0.000|x2:0.000|x3:0.652|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.348|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"

# Extend style (A column numeric label formatting) from A18 down to A19:A23
$ws.Range("A18").Copy()
$ws.Range("A19:A23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New rows 19-23: A column sequential index, B column association string
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:1.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.143|x12:0.857|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:1.000|x17:0.000"
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "x1:0.000|x2:0.000|x3:1.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "x1:0.391|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.266|x12:0.000|x13:0.343|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
